# Reverse the order of the 4 rows (players) currently in rows 6-9
# (Lauri Markkanen, Zach LaVine, John Collins, Rui Hachimura)
# so that the block becomes
# (Rui Hachimura, John Collins, Zach LaVine, Lauri Markkanen)
# while keeping their Position/Team values attached to the same player.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 6
$endRow = 9
$cols = @("A", "B", "C")

# Capture the original values for rows 6..9 before overwriting anything.
$original = @{}
for ($r = $startRow; $r -le $endRow; $r++) {
    $original[$r] = @{}
    foreach ($col in $cols) {
        $original[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Write back the rows in reverse order.
$rowCount = $endRow - $startRow + 1
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $endRow - $i
    $dstRow = $startRow + $i
    foreach ($col in $cols) {
        $ws.Range("$col$dstRow").Value2 = $original[$srcRow][$col]
    }
}
